$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23 = "InstructionHamburger": the placeholder "(       )" blank-space text
# is being replaced everywhere with an explicit TextMeshPro sprite tag for the
# hamburger-menu icon, across all language columns + the French "(F)" column.
$ws.Range("B23").Value = '(  <sprite name="FlatHamburger" ) start/pause/resume game'
$ws.Range("C23").Value = '(  <sprite name="FlatHamburger" ) Start/Pause'
$ws.Range("D23").Value = '(  <sprite name="FlatHamburger" ) comenzar/pausar/resumir juego'
$ws.Range("E23").Value = '(F) (  <sprite name="FlatHamburger" ) start/pause/resume game'
